$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column A (shifts existing A-D to B-E)
$ws.Columns("A").Insert()

# --- Header row (row 1) ---
$ws.Cells.Item(1, 1).Value = "TabName"

# --- Data row (row 2) ---
$ws.Cells.Item(2, 1).Value = "CasesTab"

$query1 = @'
MATCH (ct:clinical_trial)<--(a:arm)<--(c:case)
    WHERE c.race = "NOT_REPORTED"
WITH DISTINCT c, a, ct
RETURN 
    COALESCE(c.case_id, '') AS `Case ID`,
    COALESCE(ct.clinical_trial_designation, '') AS `Trial Code`,
    COALESCE(a.arm_id, '') AS `Arm`,
    COALESCE(a.arm_drug, '') AS `Arm Treatment`,
    COALESCE(c.disease, '') AS `Diagnosis`,
    COALESCE(c.gender, '') AS `Gender`,
    COALESCE(c.race, '') AS `Race`,
    COALESCE(c.ethnicity, '') AS `Ethnicity`
'@

$query2 = @'
MATCH (s:specimen)-->(c:case)-->(:arm)-->(ct:clinical_trial)
    WHERE c.race = "NOT_REPORTED"
OPTIONAL MATCH (f:file)-->(:sequencing_assay)-->(:nucleic_acid)-->(s)
RETURN 
    COUNT(DISTINCT f) AS number_of_files,
    COUNT(DISTINCT c.case_id) AS number_of_cases,
    COUNT(DISTINCT ct.clinical_trial_designation) AS number_of_trials
'@

$ws.Cells.Item(2, 2).Value = $query1
$ws.Cells.Item(2, 3).Value = $query2

# Wrap text on the query cells (B2, C2) to match style used before
$ws.Range("B2:C2").WrapText = $true

# Column A width - best fit narrow column
$ws.Columns("A").ColumnWidth = 7.92

# Taller row to fit the longer, wrapped query text
$ws.Rows(2).RowHeight = 174

# Selection as shown in diff
$ws.Range("B2").Select() | Out-Null
